$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.324.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4698"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2865"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06409"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.887.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7184"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.137"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "277.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.359.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007394"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.131.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.215"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.231"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.001"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.875"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.354"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09606"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.465"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.210"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.098"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04815"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.116"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6847"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.812"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.222"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.936"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4215"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8241"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.577"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.886"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "895.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05733"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.24%  "
